# Applies the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values look like plain decimal numbers (e.g. "216.92").
# The source sheet stores every Price/Volume cell as text, so for those cells we
# temporarily force a Text number format before assigning the value (otherwise Excel
# would silently convert them to numeric cells), then restore the default style so
# no stray formatting is left behind.
$forceTextCells = @("D5", "D9", "D10", "D11", "D14", "D16", "D19", "D21", "D22", "D24", "D25", "D27", "D29", "D30", "D31", "D37", "D38", "D39", "D42", "D45", "D46", "D48", "D49", "D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Cell value updates (old -> new, per diff), row by row ---
$ws.Range("D2").Value = "27.150.26"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.643.00"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "216.92"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "20.01"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.873.33"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.643.47"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "67.26"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "27.133.52"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "218.61"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "2.56"
$ws.Range("E22").Value = "  +6.06%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "147.50"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "7.53"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "15.74"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").Value = "1.263.94"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.853"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +6.58%  "
$ws.Range("D43").Value = "1.783.51"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").Value = "61.77"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "91.76"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.0975"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0976"
$ws.Range("E50").Value = "  -7.64%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +0.62%  "

# Restore default styling on the forced-text cells so no extra explicit style remains
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
